$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Cells are written as literal text (NumberFormat "@" forces text entry for
# numeric-looking strings like "66.926.91"; ClearFormats afterwards removes the
# temporary style so the cell ends up with no explicit style, matching the source.
$cellUpdates = @(
    @{ Cell = "D2"; Value = "66.926.91" }
    @{ Cell = "E2"; Value = "  -2.03%  " }
    @{ Cell = "D3"; Value = "3.480.38" }
    @{ Cell = "E3"; Value = "  -2.31%  " }
    @{ Cell = "E4"; Value = "  -0.08%  " }
    @{ Cell = "D5"; Value = "600.80" }
    @{ Cell = "E5"; Value = "  -2.90%  " }
    @{ Cell = "D6"; Value = "147.57" }
    @{ Cell = "E6"; Value = "  -4.92%  " }
    @{ Cell = "D7"; Value = "3.478.17" }
    @{ Cell = "E7"; Value = "  -2.42%  " }
    @{ Cell = "E8"; Value = "  +0.01%  " }
    @{ Cell = "E9"; Value = "  -2.40%  " }
    @{ Cell = "D10"; Value = "0.142" }
    @{ Cell = "E10"; Value = "  -2.85%  " }
    @{ Cell = "D11"; Value = "7.65" }
    @{ Cell = "E11"; Value = "  +3.35%  " }
    @{ Cell = "E12"; Value = "  -3.59%  " }
    @{ Cell = "E13"; Value = "  -4.18%  " }
    @{ Cell = "D14"; Value = "4.067.35" }
    @{ Cell = "E14"; Value = "  -2.37%  " }
    @{ Cell = "D15"; Value = "31.18" }
    @{ Cell = "E15"; Value = "  -6.20%  " }
    @{ Cell = "D16"; Value = "3.475.88" }
    @{ Cell = "E16"; Value = "  -2.73%  " }
    @{ Cell = "D17"; Value = "66.857.97" }
    @{ Cell = "E18"; Value = "  +0.27%  " }
    @{ Cell = "D19"; Value = "6.39" }
    @{ Cell = "E19"; Value = "  -5.16%  " }
    @{ Cell = "D20"; Value = "15.23" }
    @{ Cell = "E20"; Value = "  -4.88%  " }
    @{ Cell = "D21"; Value = "10.07" }
    @{ Cell = "E21"; Value = "  +0.82%  " }
    @{ Cell = "D22"; Value = "433.27" }
    @{ Cell = "E22"; Value = "  -4.61%  " }
    @{ Cell = "E23"; Value = "  -5.96%  " }
    @{ Cell = "D24"; Value = "79.07" }
    @{ Cell = "E24"; Value = "  +0.79%  " }
    @{ Cell = "E25"; Value = "  +0.11%  " }
    @{ Cell = "D26"; Value = "3.616.27" }
    @{ Cell = "E26"; Value = "  -2.45%  " }
    @{ Cell = "D27"; Value = "0.0000119" }
    @{ Cell = "E27"; Value = "  -9.02%  " }
    @{ Cell = "D28"; Value = "9.79" }
    @{ Cell = "E28"; Value = "  -7.11%  " }
    @{ Cell = "D29"; Value = "8.35" }
    @{ Cell = "E29"; Value = "  -9.09%  " }
    @{ Cell = "E30"; Value = "  -3.49%  " }
    @{ Cell = "D31"; Value = "1.57" }
    @{ Cell = "E31"; Value = "  -7.87%  " }
    @{ Cell = "E32"; Value = "  -2.72%  " }
    @{ Cell = "D34"; Value = "25.29" }
    @{ Cell = "E34"; Value = "  -3.24%  " }
    @{ Cell = "D35"; Value = "3.468.97" }
    @{ Cell = "E35"; Value = "  -2.46%  " }
    @{ Cell = "B36"; Value = "ImmutableX" }
    @{ Cell = "C36"; Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx" }
    @{ Cell = "D36"; Value = "1.79" }
    @{ Cell = "E36"; Value = "  -6.18%  " }
    @{ Cell = "B37"; Value = "NEARProtocol" }
    @{ Cell = "C37"; Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near" }
    @{ Cell = "D37"; Value = "5.90" }
    @{ Cell = "E37"; Value = "  -7.23%  " }
    @{ Cell = "E38"; Value = "  +0.02%  " }
    @{ Cell = "E39"; Value = "  -4.36%  " }
    @{ Cell = "D40"; Value = "0.999" }
    @{ Cell = "E40"; Value = "  -0.14%  " }
    @{ Cell = "D41"; Value = "173.76" }
    @{ Cell = "E41"; Value = "  -4.24%  " }
    @{ Cell = "D42"; Value = "0.0882" }
    @{ Cell = "E42"; Value = "  -3.61%  " }
    @{ Cell = "E43"; Value = "  -12.71%  " }
    @{ Cell = "E44"; Value = "  -3.79%  " }
    @{ Cell = "D45"; Value = "0.896" }
    @{ Cell = "E45"; Value = "  -0.32%  " }
    @{ Cell = "E46"; Value = "  +0.47%  " }
    @{ Cell = "D47"; Value = "28.79" }
    @{ Cell = "E47"; Value = "  -7.51%  " }
    @{ Cell = "E48"; Value = "  -6.90%  " }
    @{ Cell = "D49"; Value = "7.44" }
    @{ Cell = "E49"; Value = "  -4.55%  " }
    @{ Cell = "D50"; Value = "2.41" }
    @{ Cell = "E50"; Value = "  -9.57%  " }
    @{ Cell = "E51"; Value = "  -4.75%  " }
)

foreach ($u in $cellUpdates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.ClearFormats()
}
